# Applies the "Added BWP Bootstrap files" update: a fresh batch of Katalon
# bootstrap-run timestamps is appended to the log, and each scenario sheet's
# DateProd (column B) cell is refreshed to the latest run's timestamp.  In
# addition, sheet "VT-P-DebitCredit-DualCF-Generic" row 3's ResultProd
# (column A) flips from "Fail" to "Pass".

$wb = $excel.ActiveWorkbook

# sheet name -> hashtable of cell -> new value
$updates = @{
    "VT-P-DebitVoid-DualCF-Generic" = @{
        "B2" = "Mon Nov 10 17:28:44 IST 2025"
    }
    "VT-P-DebitVoid-SingleCF-Generic" = @{
        "B2" = "Mon Nov 10 17:34:25 IST 2025"
        "B3" = "Mon Nov 10 17:35:35 IST 2025"
        "B4" = "Mon Nov 10 17:36:42 IST 2025"
        "B5" = "Mon Nov 10 17:37:50 IST 2025"
    }
    "VT-P-DebitVoid-NoCF-Generic" = @{
        "B2" = "Mon Nov 10 17:29:50 IST 2025"
        "B3" = "Mon Nov 10 17:30:59 IST 2025"
        "B4" = "Mon Nov 10 17:32:11 IST 2025"
        "B5" = "Mon Nov 10 17:33:13 IST 2025"
    }
    "VT-P-DebitCredit-DualCF-Generic" = @{
        "A3" = "Pass"
        "B2" = "Mon Nov 10 17:15:11 IST 2025"
        "B3" = "Mon Nov 10 17:16:27 IST 2025"
        "B4" = "Mon Nov 10 17:17:30 IST 2025"
        "B5" = "Mon Nov 10 17:18:37 IST 2025"
    }
    "VT-P-DebitCredit-SingleCF-Gener" = @{
        "B2" = "Mon Nov 10 17:24:31 IST 2025"
        "B3" = "Mon Nov 10 17:25:34 IST 2025"
        "B4" = "Mon Nov 10 17:26:41 IST 2025"
        "B5" = "Mon Nov 10 17:27:44 IST 2025"
    }
    "VT-P-DebitCredit-NoCF-Generic" = @{
        "B2" = "Mon Nov 10 17:20:03 IST 2025"
        "B3" = "Mon Nov 10 17:21:08 IST 2025"
        "B4" = "Mon Nov 10 17:22:23 IST 2025"
        "B5" = "Mon Nov 10 17:23:24 IST 2025"
    }
    "VT-C-DebitCredit-DualCF-Generic" = @{
        "B2" = "Mon Nov 10 17:11:55 IST 2025"
    }
    "VT-C-DebitCredit-SingleCF-Gener" = @{
        "B2" = "Mon Nov 10 17:13:59 IST 2025"
    }
    "VT-C-DebitCredit-NoCF-Generic" = @{
        "B2" = "Mon Nov 10 17:13:00 IST 2025"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

$wb.Save()
